$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at the top. This shifts the header (previously
#    row 1) down to row 2, and all eight data rows (previously rows 2-9)
#    down to rows 3-10.
$ws.Rows.Item(1).Insert()

# 2. Per-GROUP MD5 recalculation: every data row belongs to GROUP "A", and
#    once grouped correctly all of them share the same MD5 digest that used
#    to sit on the very first data row. Unify the "MD5 & Script Version"
#    column (U) for all 8 data rows (now rows 3-10).
$unifiedMd5 = "MD5: d4f6fa9523038fdb5e8b258d4c9d18c4 | Script: v3.0.0"
for ($r = 3; $r -le 10; $r++) {
    $ws.Range("U$r").Value = $unifiedMd5
}

# 3. The MD5 column header cell lost its header style when the rows shifted
#    (it had none to begin with) -- copy the shared bold/centered/bordered
#    header formatting from a neighboring header cell so it matches the
#    rest of the header row.
$ws.Range("T2").Copy()
$ws.Range("U2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Size column U (the MD5 column) to fit its new, much wider content.
$ws.Columns.Item(21).ColumnWidth = 55.43

# 5. Update the view state: scrolled position and active selection moved
#    to the MD5 column on the header row.
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("U2").Select()
